$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.014855
$ws.Range("H2").Value = 0.044565
$ws.Range("I2").Value = 0.5536850213696453
$ws.Range("J2").Value = 0.5536850213696451
$ws.Range("M2").Value = 6.382924
$ws.Range("N2").Value = 19.148772
$ws.Range("O2").Value = 0.1363153751023214
$ws.Range("P2").Value = 0.1363153751023214
$ws.Range("Q2").Value = 0.09481833602
$ws.Range("R2").Value = 0.85336502418
$ws.Range("S2").Value = 0.07547578137654003
$ws.Range("T2").Value = 0.07547578137654003
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.014855
$ws.Range("H3").Value = 0.044565
$ws.Range("I3").Value = 0.5536850213696453
$ws.Range("J3").Value = 0.5536850213696451
$ws.Range("O3").Value = 0.6265841681043937
$ws.Range("P3").Value = 0.6265841681043938
$ws.Range("Q3").Value = 0.43583981742
$ws.Range("R3").Value = 3.92255835678
$ws.Range("S3").Value = 0.3469302685067626
$ws.Range("T3").Value = 0.3469302685067626
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.014855
$ws.Range("H4").Value = 0.044565
$ws.Range("I4").Value = 0.5536850213696453
$ws.Range("J4").Value = 0.5536850213696451
$ws.Range("O4").Value = 0.2371004567932849
$ws.Range("P4").Value = 0.2371004567932849
$ws.Range("Q4").Value = 0.1649224877666667
$ws.Range("R4").Value = 1.4843023899
$ws.Range("S4").Value = 0.1312789714863426
$ws.Range("T4").Value = 0.1312789714863426
$ws.Range("G5").Value = 0.01197433333333333
$ws.Range("I5").Value = 0.4463149786303549
$ws.Range("J5").Value = 0.4463149786303548
$ws.Range("M5").Value = 6.382924
$ws.Range("N5").Value = 19.148772
$ws.Range("O5").Value = 0.1363153751023214
$ws.Range("P5").Value = 0.1363153751023214
$ws.Range("Q5").Value = 0.07643125961733332
$ws.Range("R5").Value = 0.6878813365559999
$ws.Range("S5").Value = 0.06083959372578138
$ws.Range("T5").Value = 0.06083959372578138
$ws.Range("G6").Value = 0.01197433333333333
$ws.Range("I6").Value = 0.4463149786303549
$ws.Range("J6").Value = 0.4463149786303548
$ws.Range("O6").Value = 0.6265841681043937
$ws.Range("P6").Value = 0.6265841681043938
$ws.Range("Q6").Value = 0.3513221981639999
$ws.Range("S6").Value = 0.2796538995976312
$ws.Range("T6").Value = 0.2796538995976312
$ws.Range("G7").Value = 0.01197433333333333
$ws.Range("I7").Value = 0.4463149786303549
$ws.Range("J7").Value = 0.4463149786303548
$ws.Range("O7").Value = 0.2371004567932849
$ws.Range("P7").Value = 0.2371004567932849
$ws.Range("S7").Value = 0.1058214853069423
$ws.Range("T7").Value = 0.1058214853069423
